$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, without leaving a stray style index
# behind (values that look numeric would otherwise be auto-converted to
# real numbers by Excel's type inference).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple price (column D) updates ---
Set-TextValue $ws.Range("D2") "244.73"
Set-TextValue $ws.Range("D4") "5.196"
Set-TextValue $ws.Range("D5") "0.05730"
Set-TextValue $ws.Range("D6") "6.490"
Set-TextValue $ws.Range("D7") "3.167"
Set-TextValue $ws.Range("D9") "0.8655"
Set-TextValue $ws.Range("D10") "0.1371"
Set-TextValue $ws.Range("D11") "0.06937"
Set-TextValue $ws.Range("D12") "0.03181"
Set-TextValue $ws.Range("D13") "0.02931"
Set-TextValue $ws.Range("D14") "0.09323"
Set-TextValue $ws.Range("D15") "3.833"
Set-TextValue $ws.Range("D16") "0.001529"
Set-TextValue $ws.Range("D17") "0.04717"

# --- Rows 18-24: coin list shifted up by one, with "One" re-inserted at row 18 ---
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0005981"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D19") "0.006149"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D20") "0.001238"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D21") "0.004109"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D22") "0.00008503"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D23") "3.550"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D24") "2.161"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# --- More simple price updates ---
Set-TextValue $ws.Range("D25") "0.3193"
Set-TextValue $ws.Range("D27") "0.0002330"

# --- Rows 41-43: coin list shifted, with rotation among KickToken/BKEXToken/CEJI ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1053"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.002223"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "0.003066"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- More simple updates ---
Set-TextValue $ws.Range("D44") "0.008117"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

Set-TextValue $ws.Range("D45") "0.00005455"

Set-TextValue $ws.Range("D47") "0.4539"
Set-TextValue $ws.Range("D48") "0.002566"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("D50") "0.0002000"
